$wb = $excel.ActiveWorkbook

$loginWs = $wb.Worksheets.Item("Login")
$loginWs.Range("G2").Value = "Success - 2020/12/19 16:28:07"
$loginWs.Range("G3").Value = "Success - 2020/12/19 16:28:11"

$schoolWs = $wb.Worksheets.Item("School Search")
$schoolWs.Range("C2").Value = "Success - 2020/12/19 16:28:18"
$schoolWs.Range("C3").Value = "Success - 2020/12/19 16:28:21"

$productWs = $wb.Worksheets.Item("Product Search")
$productWs.Range("K2").Value = "Success - 2020/12/19 16:28:48"
$productWs.Range("K3").Value = "Success - 2020/12/19 16:29:10"
$productWs.Range("K4").Value = "Success - 2020/12/19 16:29:33"

$cartWs = $wb.Worksheets.Item("Shopping Cart")
$cartWs.Range("G2").Value = "Success - 2020/12/19 16:29:35"
$cartWs.Range("G3").Value = "Success - 2020/12/19 16:29:35"
$cartWs.Range("G4").Value = "Success - 2020/12/19 16:29:36"

$checkoutWs = $wb.Worksheets.Item("Checkout")
$checkoutWs.Range("P2").Value = "Success - 2020/12/19 16:29:46"
$checkoutWs.Range("P3").Value = "Success - 2020/12/19 16:29:53"
$checkoutWs.Range("P4").Value = "Success - 2020/12/19 16:30:02"

$paymentWs = $wb.Worksheets.Item("Payment")
# This cell previously carried a quote-prefix (text-forced) style; prefix the
# new value with a literal apostrophe so the engine keeps treating it as an
# explicitly-quoted text entry and preserves that cell style.
$paymentWs.Range("C2").Value = "'Success - 2020/12/19 16:30:12"
